$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.238.82"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "3.884.72"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "472.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.743"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000316"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.68%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.526.23"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.20%  "
$ws.Range("D16").Value = "3.876.94"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +6.59%  "
$ws.Range("D20").Value = "67.522.85"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "729.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  +8.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.94%  "
$ws.Range("E35").Value = "  +7.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0483"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.350"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.34%  "
$ws.Range("D41").Value = "0.0₃0690"
$ws.Range("E41").Value = "  -7.78%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.55%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.51%  "
